# Updated symbol list on Wed Dec 28 15:49:31 UTC 2022 with GitHub Actions
#
# Prices in column D are stored as TEXT (not numbers), so every new price is
# written with a leading apostrophe to force Excel to keep it as text
# (this mirrors how the source workbook stores e.g. "243.34" as a string).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- rows 2-9: price (D) refreshed, everything else unchanged ---
$ws.Range("D2").Formula  = "'243.78"
$ws.Range("D3").Formula  = "'23.80"
$ws.Range("D4").Formula  = "'5.272"
$ws.Range("D5").Formula  = "'0.05846"
$ws.Range("D6").Formula  = "'6.481"
$ws.Range("D7").Formula  = "'3.343"
$ws.Range("D8").Formula  = "'0.8111"
$ws.Range("D9").Formula  = "'0.8917"

# --- rows 10-18: "One" jumps from row 18 up to row 10; WazirX..CoinExToken
#     each shift down by one row (keeping name/link) and get a refreshed
#     price + rank-prefixed volume label ---
$ws.Range("B10").Formula = "One"
$ws.Range("C10").Formula = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D10").Formula = "'0.01035"
$ws.Range("E10").Formula = "9OneONE"

$ws.Range("B11").Formula = "WazirX"
$ws.Range("C11").Formula = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").Formula = "'0.1389"
$ws.Range("E11").Formula = "10WazirXWRX"

$ws.Range("B12").Formula = "MandalaExchangeToken"
$ws.Range("C12").Formula = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Formula = "'0.07235"
$ws.Range("E12").Formula = "11MandalaExchangeTokenMDX"

$ws.Range("B13").Formula = "LiechtensteinCryptoassetsExchange"
$ws.Range("C13").Formula = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D13").Formula = "'0.03097"
$ws.Range("E13").Formula = "12LiechtensteinCryptoassetsExchangeLCX"

$ws.Range("B14").Formula = "BitrueCoin"
$ws.Range("C14").Formula = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D14").Formula = "'0.03021"
$ws.Range("E14").Formula = "13BitrueCoinBTR"

$ws.Range("B15").Formula = "BitMartToken"
$ws.Range("C15").Formula = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D15").Formula = "'0.09324"
$ws.Range("E15").Formula = "14BitMartTokenBMX"

$ws.Range("B16").Formula = "MCDex"
$ws.Range("C16").Formula = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D16").Formula = "'3.836"
$ws.Range("E16").Formula = "15MCDexMCB"

$ws.Range("B17").Formula = "BitForexToken"
$ws.Range("C17").Formula = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D17").Formula = "'0.001569"
$ws.Range("E17").Formula = "16BitForexTokenBF"

$ws.Range("B18").Formula = "CoinExToken"
$ws.Range("C18").Formula = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D18").Formula = "'0.04728"
$ws.Range("E18").Formula = "17CoinExTokenCET"

# --- remaining rows: only the price (D) was refreshed ---
$ws.Range("D19").Formula = "'0.006226"
$ws.Range("D20").Formula = "'0.001262"
$ws.Range("D21").Formula = "'0.003852"
$ws.Range("D22").Formula = "'0.00008715"
$ws.Range("D23").Formula = "'3.552"
$ws.Range("D24").Formula = "'2.177"
$ws.Range("D25").Formula = "'0.3196"
$ws.Range("D26").Formula = "'0.1314"
$ws.Range("D28").Formula = "'0.0002344"
$ws.Range("D40").Formula = "'0.03796"
$ws.Range("D41").Formula = "'0.006351"
$ws.Range("D42").Formula = "'0.1052"
$ws.Range("D43").Formula = "'0.002522"
$ws.Range("D44").Formula = "'0.007092"
$ws.Range("D45").Formula = "'0.00005352"
$ws.Range("D47").Formula = "'0.5510"
$ws.Range("D48").Formula = "'0.01668"
$ws.Range("D49").Formula = "'0.00002104"
$ws.Range("D50").Formula = "'0.0002004"
